# PAS-2715: Update VIN Upload CA_SELECT file
#   - Drop the STAT column
#   - Drop the CHOICE_TIER column
#   - Append eight new columns: BI_SYMBOL, PD_SYMBOL, UM_SYMBOL, MP_SYMBOL,
#     ENTRYDATE, VALID, ANTITHEFT_DISCOUNT, RESTRAINTS_DISCOUNT

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the STAT column (Z) - everything to its right shifts left one column
$ws.Columns("Z").Delete()

# Remove the CHOICE_TIER column (now at AB after the previous delete)
$ws.Columns("AB").Delete()

# ---- Append the new trailing columns (headers in row 1, sample row in row 2) ----

# Headers: copy the header style from the last existing header cell (AB1)
$ws.Range("AC1:AJ1").Style = $ws.Range("AB1").Style
$ws.Range("AC1").Value = "BI_SYMBOL"
$ws.Range("AD1").Value = "PD_SYMBOL"
$ws.Range("AE1").Value = "UM_SYMBOL"
$ws.Range("AF1").Value = "MP_SYMBOL"
$ws.Range("AG1").Value = "ENTRYDATE"
$ws.Range("AH1").Value = "VALID"
$ws.Range("AI1").Value = "ANTITHEFT_DISCOUNT"
$ws.Range("AJ1").Value = "RESTRAINTS_DISCOUNT"

# Data row: AC2:AF2 share the style used by the other data cells (copy from AB2)
$ws.Range("AC2:AF2").Style = $ws.Range("AB2").Style
$ws.Range("AC2").Value = "K"
$ws.Range("AD2").Value = "K"
$ws.Range("AE2").Value = "K"
$ws.Range("AF2").Value = "K"
$ws.Range("AG2").Value = 20000101
$ws.Range("AH2").Value = "Y"
$ws.Range("AI2").Value = "Y"
$ws.Range("AJ2").Value = "N"

# Widen the new RESTRAINTS_DISCOUNT column like the source workbook
$ws.Columns("AJ").ColumnWidth = 20.6

# Scroll the view over to the new columns and leave the cursor on AJ12,
# matching the author's saved view state
$ws.Range("AJ12").Select()
$excel.ActiveWindow.ScrollColumn = 25
